$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.167.06"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.855.83"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +1.19%  "
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "'310.75"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.4782"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").Value = "'0.3707"
$ws.Range("D9").Value = "'0.07285"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'0.9350"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "'19.98"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'0.07826"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "1.876.48"
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").Value = "'5.404"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "'6.515"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "'89.77"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "'1.019"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'0.000008716"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "27.172.44"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "'14.66"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").Value = "'5.083"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'1.939"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "'153.29"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").Value = "'1.995"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").Value = "'115.31"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'4.932"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "'0.08873"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'3.306"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("D33").Value = "'4.554"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'0.7372"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "'2.691"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("E36").Value = "  +2.85%  "
$ws.Range("D37").Value = "'0.01999"
$ws.Range("E37").Value = "  +2.85%  "
$ws.Range("D38").Value = "'0.05248"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").Value = "'0.5311"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").Value = "'7.054"
$ws.Range("E40").Value = "  +2.07%  "
$ws.Range("D41").Value = "'0.1527"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("D42").Value = "'8.329"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").Value = "'10.65"
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("D44").Value = "'0.4766"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").Value = "'1.018"
$ws.Range("D46").Value = "'102.22"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "'1.629"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "'66.06"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("D49").Value = "'0.06068"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "'0.8956"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'36.59"
$ws.Range("E51").Value = "  +0.90%  "
